$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose District (column G) value is standardized to the official
# name "Belagavi (Belgaum)" (various misspellings: Belgaum, Belagavi,
# Belgavi, and the one stray "Chikkodi" in row 37 that was corrected too).
$rowsToStandardize = @(
    4, 5, 6, 7, 9, 10, 11, 12, 13, 14, 15, 17,
    19, 20, 22, 23, 24, 25, 26, 27, 28, 29, 30,
    31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41,
    42, 43, 44, 45, 46, 48, 49, 50, 51, 52, 53,
    54, 55, 56
)

foreach ($r in $rowsToStandardize) {
    $ws.Range("G$r").Value = "Belagavi (Belgaum)"
}

# Row 18 had a stray empty F-column cell (no real address data); remove it.
$ws.Range("F18").ClearContents()
